$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 4
$ws.Range("F8").Value = 0
$ws.Range("F14").Value = 4
$ws.Range("F18").Value = 1
$ws.Range("F20").Value = -1
$ws.Range("F31").Value = 2
$ws.Range("F32").Value = 1
$ws.Range("F38").Value = 3
$ws.Range("F40").Value = 0
$ws.Range("F42").Value = 1
$ws.Range("F54").Value = -9
$ws.Range("F57").Value = 5
